$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 -- shifts existing rows 10..56 down to 11..57,
# carrying their values/formatting with them (matches dimension A1:R56 -> A1:R57).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted (blank) row 10 with the new weekly record.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 44847
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 100112026
$ws.Range("G10").Value = "Haba"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 9000
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 9000
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 360
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
